$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C is entirely empty; delete it so column D (B-S65-20_k) shifts
# into C and column E (B-S65-20_m) shifts into D.
$ws.Range("C:C").Delete() | Out-Null

# Update the active selection to match the target state.
$ws.Range("M9").Select() | Out-Null
